$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the existing "Add Panels" sheet, placing the copy right after
#    it, then rename the copy to "Sheet1". This preserves the original,
#    un-edited row 8 data (2.2824E-2 / 2.7389E-2 values) on a second tab.
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("Add Panels")
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "Sheet1"
$newSheet.Range("D5").Select()

# ---------------------------------------------------------------------------
# 2. Re-activate the original "Add Panels" sheet and make the edits to it.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Add Panels")
$ws.Activate()

# Update the existing FIRECLASS 64-2 row (row 8): Stand By Hours / Battery
# Factor columns move from fractional values to plain Ah numbers.
$ws.Range("F8").Value = 22.9
$ws.Range("C5").Copy()
$ws.Range("F8").PasteSpecial(-4122)

$ws.Range("J8").Value = 28.98
$ws.Range("C5").Copy()
$ws.Range("J8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Add the new FC702S row (row 9), mirroring row 8's layout/styling.
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "FC702S"
$ws.Range("B9").Value = "Node1"
$ws.Range("D9").Value = "FIM"
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 25.12
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = 77
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "1.75"
$ws.Range("J9").Value = 37.58
$ws.Range("K9").Value = "Minimum Battery size(Ah)"

# Copy each row-8 cell's formatting down onto the matching row-9 cell.
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("E8").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("G8").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("H8").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("I8").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("J8").Copy()
$ws.Range("J9").PasteSpecial(-4122)
$ws.Range("K8").Copy()
$ws.Range("K9").PasteSpecial(-4122)

$ws.Range("J9").Select()
